# issue #5: stock data from json to db
#
# The "股票" (stock) sheet is extended with three new metadata columns
# that the json->db importer now emits: "category" (inserted right after
# "property_category"), and "source_file" / "index" appended at the end.
# Also fixes a stray leading "、" on one of the existing text values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

# Insert a new column before the old "date" column (I) so the layout
# becomes: name, owner, quantity, face_value, currency, total,
# property_category, category, date, legislator_name, legislator_id,
# source_file, index
$ws.Columns.Item(9).Insert()

# Header row
$ws.Range("I1").Value = "category"
$ws.Range("M1").Value = "source_file"
$ws.Range("N1").Value = "index"

# Match the header formatting (bold + border) used by the rest of row 1.
$ws.Range("L1").Copy()
$ws.Range("M1:N1").PasteSpecial(-4122)

# Data rows: the sheet has 7 data rows (rows 2-8). Column A holds the
# original row index (82-88), which is reused for the new "index" column.
for ($r = 2; $r -le 8; $r++) {
    $ws.Range("I" + $r).Value = "normal"
    $ws.Range("M" + $r).Value = "tmp30a51"
    $idx = $ws.Range("A" + $r).Value()
    $ws.Range("N" + $r).Value = $idx
}

# Fix the stray leading "、" character on G7 ("、2100000" -> "2100000"),
# keeping the cell as text (it was stored as text before the fix too).
$ws.Range("G7").NumberFormat = "@"
$ws.Range("G7").Value = "2100000"
